$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values per row (case with 380 kV), row number => column => value
$data = @{
    2 = @{ B=1.02; C=1.043628789963518; D=1.044186611904769; E=1.056918513702797; F=1.06390786281133; I=1.037222416452245; J=1.048698214599805; K=1.046958369591109; L=1.059654907143038; M=1.066625232154586; N=1.050187486351903 }
    3 = @{ B=1.02; C=1.04528104831694; D=1.045438760760032; E=1.058323288267223; F=1.06541934479553; I=1.037608732202804; J=1.049993916369561; K=1.04802055457833; L=1.060871913041368; M=1.067950071087328; N=1.051485028166791 }
    4 = @{ B=1.02; C=1.046347741975253; D=1.04624666998957; E=1.059231293567985; F=1.066396169807336; I=1.037856337250241; J=1.050829583227988; K=1.048704995129813; L=1.061657894302469; M=1.068805620376381; N=1.05232188176801 }
    5 = @{ B=1.02; C=1.046795610459589; D=1.0465857687135; E=1.05961279214777; F=1.066806546859658; I=1.037959866773084; J=1.051180250882046; K=1.048992055591163; L=1.061987967712338; M=1.06916489150131; N=1.052673047410359 }
    6 = @{ B=1.02; C=1.046870776490634; D=1.046642673072088; E=1.059676834344699; F=1.06687543478174; I=1.037977216871366; J=1.05123909186158; K=1.049040214739288; L=1.062043367977148; M=1.069225191362398; N=1.052731971950832 }
    7 = @{ B=1.02; C=1.046353728637518; D=1.046251203179027; E=1.059236392051909; F=1.066401654374971; I=1.037857722827242; J=1.050834271394941; K=1.048708833500975; L=1.061662306137649; M=1.068810422542906; N=1.052326576592698 }
    8 = @{ B=1.02; C=1.044187688226992; D=1.044610265515618; E=1.057393470667892; F=1.064418929704601; I=1.037353465627501; J=1.049136675602168; K=1.047317937459875; L=1.060066515219684; M=1.067073326352768; N=1.050626570019182 }
    9 = @{ B=1.02; C=1.04035172465163; D=1.04170064098358; E=1.054138193448349; F=1.060915514917693; I=1.036446633657535; J=1.04612391386863; K=1.044844733994035; L=1.057242715638188; M=1.063998901285798; N=1.047609529818414 }
    10 = @{ B=1.02; C=1.037780816354047; D=1.039748235153159; E=1.05196229666549; F=1.058572926461472; I=1.035829603902422; J=1.044100461375183; K=1.043180494962868; L=1.055351816769756; M=1.06193976696596; N=1.045583203790332 }
    11 = @{ B=1.02; C=1.036664193236477; D=1.038899708149144; E=1.05101865180564; F=1.057556787198023; I=1.035559420615947; J=1.043220611859421; K=1.042456095186066; L=1.054530957314916; M=1.061045773647551; N=1.044702104787338 }
    12 = @{ B=1.02; C=1.036248904070339; D=1.038584048592972; E=1.050667910726249; F=1.057179069757234; I=1.035458607163529; J=1.042893232418335; K=1.04218644485152; L=1.05422573162288; M=1.060713338578848; N=1.044374260429895 }
    13 = @{ B=1.02; C=1.036338009024149; D=1.038651780465593; E=1.050743156401881; F=1.057260104178734; I=1.035480252658812; J=1.042963482128274; K=1.042244311976255; L=1.054291218272239; M=1.060784663750185; N=1.04444460990248 }
    14 = @{ B=1.02; C=1.036629876110701; D=1.038873625456799; E=1.050989664169513; F=1.057525570688526; I=1.035551096655843; J=1.043193562145869; K=1.042433817608628; L=1.054505733879683; M=1.061018301984982; N=1.04467501666009 }
    15 = @{ B=1.02; C=1.036809634972777; D=1.039010247732735; E=1.051141515079452; F=1.057689096212327; I=1.035594685526902; J=1.043335246948447; K=1.042550501683015; L=1.054637861115917; M=1.061162205412121; N=1.04481690267134 }
    16 = @{ B=1.02; C=1.03785485003422; D=1.039804482481406; E=1.05202489159436; F=1.058640325843706; I=1.035847471432233; J=1.044158775586286; K=1.043228490716019; L=1.0554062498167; M=1.061999047529095; N=1.045641600814304 }
    17 = @{ B=1.02; C=1.038509564677388; D=1.040301842156377; E=1.05257861125379; F=1.059236522530102; I=1.03600523004666; J=1.044674359977269; K=1.043652758705302; L=1.055887674973984; M=1.062523334137406; N=1.046157917394261 }
    18 = @{ B=1.02; C=1.038891121304108; D=1.040591643037928; E=1.052901445608273; F=1.059584102630851; I=1.036096958267842; J=1.044974737199197; K=1.043899863597761; L=1.056168281247644; M=1.062828913173848; N=1.046458721186298 }
    19 = @{ B=1.02; C=1.039021167162497; D=1.040690406955043; E=1.05301150014076; F=1.059702589707619; I=1.036128186196256; J=1.045077098285613; K=1.043984058607271; L=1.056263926912521; M=1.062933069258186; N=1.046561227637197 }
    20 = @{ B=1.02; C=1.038439353955384; D=1.040248511333351; E=1.052519217018943; F=1.059172574064248; I=1.035988334034082; J=1.044619079392132; K=1.043607276424231; L=1.055836043448174; M=1.062467106846016; N=1.046102558304352 }
    21 = @{ B=1.02; C=1.036543943154268; D=1.03880831093347; E=1.050917080179689; F=1.057447405176588; I=1.035530247465331; J=1.043125824978375; K=1.042378028897339; L=1.054442573333998; M=1.060949511531402; N=1.044607183298049 }
    22 = @{ B=1.02; C=1.035349174837594; D=1.037900024424536; E=1.049908421651782; F=1.056361109047541; I=1.03523959381533; J=1.042183687054715; K=1.041601815028599; L=1.05356457655262; M=1.059993216001792; N=1.043663707430473 }
    23 = @{ B=1.02; C=1.035982837911154; D=1.038381790638369; E=1.050443259837543; F=1.056937131595188; I=1.035393926009059; J=1.042683445846802; K=1.042013620076422; L=1.054030198995039; M=1.060500370762317; N=1.044164175937367 }
    24 = @{ B=1.02; C=1.038471080169916; D=1.040272610165599; E=1.052546055152785; F=1.05920147015291; I=1.035995969510105; J=1.044644059404217; K=1.043627829032358; L=1.055859374133061; M=1.062492514244803; N=1.04612757379092 }
    25 = @{ B=1.02; C=1.041345755358387; D=1.042455042558827; E=1.054980733578801; F=1.06182242296154; I=1.036683255943923; J=1.046905377023637; K=1.045486802772753; L=1.05797417872744; M=1.064795355829644; N=1.048392102740736 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
